$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target cells are stored as text (e.g. "9.00"), not numbers.
# Force text number format first so Excel keeps the value as a string
# instead of converting it to a numeric value.
$targets = @{
    "B2" = "9.00";  "D2" = "9.00";
    "B3" = "28.00"; "D3" = "28.00";
    "B4" = "14.00"; "D4" = "14.00";
    "B5" = "9.00";  "D5" = "9.00";
    "B6" = "22.00"; "D6" = "22.00";
    "B7" = "82.00"; "D7" = "82.00";
}

foreach ($addr in $targets.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $targets[$addr]
}
